# Refresh the "Price" (D) and "Volume(1h)" (E) columns on Sheet1 with the
# latest scraped crypto values from this GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.045.20'
$ws.Range("E2").Value = '  -1.63%  '

$ws.Range("D3").Value = '3.743.94'
$ws.Range("E3").Value = '  +0.85%  '

$ws.Range("E4").Value = '  +0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '623.88'
$ws.Range("E5").Value = '  +0.39%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.71'
$ws.Range("E6").Value = '  -1.49%  '

$ws.Range("D7").Value = '3.740.92'
$ws.Range("E7").Value = '  +0.97%  '

$ws.Range("E8").Value = '  +0.05%  '

$ws.Range("E9").Value = '  -1.32%  '

$ws.Range("E10").Value = '  +1.81%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.30'
$ws.Range("E11").Value = '  -5.28%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.487'
$ws.Range("E12").Value = '  -3.28%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '41.05'
$ws.Range("E13").Value = '  +1.07%  '

$ws.Range("E14").Value = '  +1.24%  '

$ws.Range("D15").Value = '4.363.44'
$ws.Range("E15").Value = '  +1.38%  '

$ws.Range("D16").Value = '3.737.67'
$ws.Range("E16").Value = '  +1.15%  '

$ws.Range("D17").Value = '70.076.97'
$ws.Range("E17").Value = '  -1.54%  '

$ws.Range("E18").Value = '  -1.23%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.59'
$ws.Range("E19").Value = '  +0.30%  '

$ws.Range("E20").Value = '  -0.91%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '506.74'
$ws.Range("E21").Value = '  -2.86%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.44'
$ws.Range("E22").Value = '  +0.69%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.725'
$ws.Range("E23").Value = '  -2.71%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.53'
$ws.Range("E24").Value = '  -0.16%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.76'
$ws.Range("E25").Value = '  -2.21%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.12'
$ws.Range("E26").Value = '  -3.64%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.24'
$ws.Range("E27").Value = '  +0.67%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000136'
$ws.Range("E28").Value = '  +21.67%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.06%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.49'
$ws.Range("E30").Value = '  -2.30%  '

$ws.Range("E31").Value = '  +0.56%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.91'
$ws.Range("E32").Value = '  -3.39%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.25'
$ws.Range("E33").Value = '  -2.06%  '

$ws.Range("E34").Value = '  -0.77%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.13%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.07'
$ws.Range("E36").Value = '  +2.40%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.22'
$ws.Range("E37").Value = '  +0.93%  '

$ws.Range("E38").Value = '  +1.48%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.337'
$ws.Range("E39").Value = '  -3.51%  '

$ws.Range("E40").Value = '  -6.81%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '50.48'
$ws.Range("E41").Value = '  -2.26%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '45.08'
$ws.Range("E42").Value = '  -0.62%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '425.45'
$ws.Range("E43").Value = '  -2.61%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.72'
$ws.Range("E44").Value = '  -1.48%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.81'
$ws.Range("E45").Value = '  -1.27%  '

$ws.Range("D46").Value = '3.003.13'
$ws.Range("E46").Value = '  -4.34%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0364'
$ws.Range("E47").Value = '  -1.53%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.36'
$ws.Range("E48").Value = '  -3.81%  '

$ws.Range("E49").Value = '  -0.04%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '137.24'
$ws.Range("E50").Value = '  -2.25%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.53'
$ws.Range("E51").Value = '  +1.75%  '
